# Update NATMI LR-pair (Col1a1-Itga11) TPM-derived statistics with newly
# recomputed TPM values. Only the numeric result columns (G..T) for data
# rows 2..17 change; identifier columns A..F stay the same.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 25.35940266666667
$ws.Range("H2").Value = 76.078208
$ws.Range("I2").Value = 0.005186643687654987
$ws.Range("J2").Value = 0.005186643687654986
$ws.Range("M2").Value = 0.152959
$ws.Range("N2").Value = 0.458877
$ws.Range("O2").Value = 0.004761500378002596
$ws.Range("P2").Value = 0.004761500378002596
$ws.Range("Q2").Value = 3.878948872490667
$ws.Range("R2").Value = 34.910539852416
$ws.Range("S2").Value = 0.000024696205879334
$ws.Range("T2").Value = 0.000024696205879334
$ws.Range("G3").Value = 25.35940266666667
$ws.Range("H3").Value = 76.078208
$ws.Range("I3").Value = 0.005186643687654987
$ws.Range("J3").Value = 0.005186643687654986
$ws.Range("O3").Value = 0.9837878817404418
$ws.Range("P3").Value = 0.9837878817404418
$ws.Range("Q3").Value = 801.4412667647148
$ws.Range("R3").Value = 7212.971400882433
$ws.Range("S3").Value = 0.005102557206820534
$ws.Range("T3").Value = 0.005102557206820533
$ws.Range("G4").Value = 25.35940266666667
$ws.Range("H4").Value = 76.078208
$ws.Range("I4").Value = 0.005186643687654987
$ws.Range("J4").Value = 0.005186643687654986
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2671263333333334
$ws.Range("N4").Value = 0.8013790000000001
$ws.Range("O4").Value = 0.008315444904458803
$ws.Range("P4").Value = 0.008315444904458805
$ws.Range("Q4").Value = 6.774164249870223
$ws.Range("R4").Value = 60.96747824883201
$ws.Range("S4").Value = 0.00004312924982375408
$ws.Range("T4").Value = 0.00004312924982375408
$ws.Range("G5").Value = 25.35940266666667
$ws.Range("H5").Value = 76.078208
$ws.Range("I5").Value = 0.005186643687654987
$ws.Range("J5").Value = 0.005186643687654986
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1007146666666667
$ws.Range("N5").Value = 0.302144
$ws.Range("O5").Value = 0.00313517297709673
$ws.Range("P5").Value = 0.00313517297709673
$ws.Range("Q5").Value = 2.554063786439111
$ws.Range("R5").Value = 22.986574077952
$ws.Range("S5").Value = 0.00001626102513136525
$ws.Range("T5").Value = 0.00001626102513136525
$ws.Range("I6").Value = 0.9837462940761621
$ws.Range("J6").Value = 0.983746294076162
$ws.Range("M6").Value = 0.152959
$ws.Range("N6").Value = 0.458877
$ws.Range("O6").Value = 0.004761500378002596
$ws.Range("P6").Value = 0.004761500378002596
$ws.Range("Q6").Value = 735.7169314148255
$ws.Range("R6").Value = 6621.452382733429
$ws.Range("S6").Value = 0.004684108351102299
$ws.Range("T6").Value = 0.004684108351102298
$ws.Range("I7").Value = 0.9837462940761621
$ws.Range("J7").Value = 0.983746294076162
$ws.Range("O7").Value = 0.9837878817404418
$ws.Range("P7").Value = 0.9837878817404418
$ws.Range("S7").Value = 0.9677976828191973
$ws.Range("T7").Value = 0.9677976828191972
$ws.Range("I8").Value = 0.9837462940761621
$ws.Range("J8").Value = 0.983746294076162
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2671263333333334
$ws.Range("N8").Value = 0.8013790000000001
$ws.Range("O8").Value = 0.008315444904458803
$ws.Range("P8").Value = 0.008315444904458805
$ws.Range("Q8").Value = 1284.849968031262
$ws.Range("R8").Value = 11563.64971228136
$ws.Range("S8").Value = 0.008180288108355853
$ws.Range("T8").Value = 0.008180288108355855
$ws.Range("I9").Value = 0.9837462940761621
$ws.Range("J9").Value = 0.983746294076162
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1007146666666667
$ws.Range("N9").Value = 0.302144
$ws.Range("O9").Value = 0.00313517297709673
$ws.Range("P9").Value = 0.00313517297709673
$ws.Range("Q9").Value = 484.4271047043129
$ws.Range("R9").Value = 4359.843942338815
$ws.Range("S9").Value = 0.003084214797506636
$ws.Range("T9").Value = 0.003084214797506636
$ws.Range("G10").Value = 51.27300266666666
$ws.Range("H10").Value = 153.819008
$ws.Range("I10").Value = 0.01048663484403512
$ws.Range("J10").Value = 0.01048663484403512
$ws.Range("M10").Value = 0.152959
$ws.Range("N10").Value = 0.458877
$ws.Range("O10").Value = 0.004761500378002596
$ws.Range("P10").Value = 0.004761500378002596
$ws.Range("Q10").Value = 7.842667214890667
$ws.Range("R10").Value = 70.58400493401601
$ws.Range("S10").Value = 0.00004993211577384844
$ws.Range("T10").Value = 0.00004993211577384844
$ws.Range("G11").Value = 51.27300266666666
$ws.Range("H11").Value = 153.819008
$ws.Range("I11").Value = 0.01048663484403512
$ws.Range("J11").Value = 0.01048663484403512
$ws.Range("O11").Value = 0.9837878817404418
$ws.Range("P11").Value = 0.9837878817404418
$ws.Range("Q11").Value = 1620.397008089515
$ws.Range("R11").Value = 14583.57307280563
$ws.Range("S11").Value = 0.01031662427979882
$ws.Range("T11").Value = 0.01031662427979882
$ws.Range("G12").Value = 51.27300266666666
$ws.Range("H12").Value = 153.819008
$ws.Range("I12").Value = 0.01048663484403512
$ws.Range("J12").Value = 0.01048663484403512
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.2671263333333334
$ws.Range("N12").Value = 0.8013790000000001
$ws.Range("O12").Value = 0.008315444904458803
$ws.Range("P12").Value = 0.008315444904458805
$ws.Range("Q12").Value = 13.69636920133689
$ws.Range("R12").Value = 123.267322812032
$ws.Range("S12").Value = 0.00008720103427875202
$ws.Range("T12").Value = 0.00008720103427875203
$ws.Range("G13").Value = 51.27300266666666
$ws.Range("H13").Value = 153.819008
$ws.Range("I13").Value = 0.01048663484403512
$ws.Range("J13").Value = 0.01048663484403512
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1007146666666667
$ws.Range("N13").Value = 0.302144
$ws.Range("O13").Value = 0.00313517297709673
$ws.Range("P13").Value = 0.00313517297709673
$ws.Range("Q13").Value = 5.163943372572444
$ws.Range("R13").Value = 46.475490353152
$ws.Range("S13").Value = 0.0000328774141836999
$ws.Range("T13").Value = 0.00003287741418369991
$ws.Range("G14").Value = 2.837922333333333
$ws.Range("H14").Value = 8.513767
$ws.Range("I14").Value = 0.0005804273921477663
$ws.Range("J14").Value = 0.0005804273921477662
$ws.Range("M14").Value = 0.152959
$ws.Range("N14").Value = 0.458877
$ws.Range("O14").Value = 0.004761500378002596
$ws.Range("P14").Value = 0.004761500378002596
$ws.Range("Q14").Value = 0.4340857621843334
$ws.Range("R14").Value = 3.906771859659
$ws.Range("S14").Value = 0.000002763705247114651
$ws.Range("T14").Value = 0.00000276370524711465
$ws.Range("G15").Value = 2.837922333333333
$ws.Range("H15").Value = 8.513767
$ws.Range("I15").Value = 0.0005804273921477663
$ws.Range("J15").Value = 0.0005804273921477662
$ws.Range("O15").Value = 0.9837878817404418
$ws.Range("P15").Value = 0.9837878817404418
$ws.Range("Q15").Value = 89.68776195963534
$ws.Range("R15").Value = 807.189857636718
$ws.Range("S15").Value = 0.0005710174346251798
$ws.Range("T15").Value = 0.0005710174346251796
$ws.Range("G16").Value = 2.837922333333333
$ws.Range("H16").Value = 8.513767
$ws.Range("I16").Value = 0.0005804273921477663
$ws.Range("J16").Value = 0.0005804273921477662
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2671263333333334
$ws.Range("N16").Value = 0.8013790000000001
$ws.Range("O16").Value = 0.008315444904458803
$ws.Range("P16").Value = 0.008315444904458805
$ws.Range("Q16").Value = 0.7580837871881111
$ws.Range("R16").Value = 6.822754084693
$ws.Range("S16").Value = 0.000004826512000443455
$ws.Range("T16").Value = 0.000004826512000443455
$ws.Range("G17").Value = 2.837922333333333
$ws.Range("H17").Value = 8.513767
$ws.Range("I17").Value = 0.0005804273921477663
$ws.Range("J17").Value = 0.0005804273921477662
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1007146666666667
$ws.Range("N17").Value = 0.302144
$ws.Range("O17").Value = 0.00313517297709673
$ws.Range("P17").Value = 0.00313517297709673
$ws.Range("Q17").Value = 0.2858204018275555
$ws.Range("R17").Value = 2.572383616448
$ws.Range("S17").Value = 0.000001819740275028404
$ws.Range("T17").Value = 0.000001819740275028403
